$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the current row 203 so the old row 203 data
# (now moved) ends up at row 205, and two new rows of data land at 203/204.
$ws.Rows.Item(204).Insert()
$ws.Rows.Item(204).Insert()

# Row 203: update in place (A/B/C already correct) with new data
$ws.Range("D203").Value = 44595
$ws.Range("E203").Value = 10
$ws.Range("F203").Value = 100112024
$ws.Range("G203").Value = "Choclo"
$ws.Range("H203").Value = "Choclero"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 3000
$ws.Range("K203").Value = 300
$ws.Range("L203").Value = 300
$ws.Range("M203").Value = 300
$ws.Range("N203").Value = "$/unidad"
$ws.Range("O203").Value = "Región del Maule"
$ws.Range("P203").Value = 300
$ws.Range("Q203").Value = 1
$ws.Range("R203").Value = "Hortaliza"

# Row 204: new row
$ws.Range("A204").Value = 4
$ws.Range("B204").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C204").Value = "Los Lagos"
$ws.Range("D204").Value = 44595
$ws.Range("E204").Value = 10
$ws.Range("F204").Value = 100112024
$ws.Range("G204").Value = "Choclo"
$ws.Range("H204").Value = "Dulce o Americano"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 3000
$ws.Range("K204").Value = 250
$ws.Range("L204").Value = 250
$ws.Range("M204").Value = 250
$ws.Range("N204").Value = "$/unidad"
$ws.Range("O204").Value = "Región de O'Higgins"
$ws.Range("P204").Value = 250
$ws.Range("Q204").Value = 1
$ws.Range("R204").Value = "Hortaliza"

# Row 205: this is the original row 203 data (unchanged), now shifted down
$ws.Range("A205").Value = 4
$ws.Range("B205").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C205").Value = "Los Lagos"
$ws.Range("D205").Value = 44544
$ws.Range("E205").Value = 10
$ws.Range("F205").Value = 100112024
$ws.Range("G205").Value = "Choclo"
$ws.Range("H205").Value = "Dulce o Americano"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 200
$ws.Range("K205").Value = 24000
$ws.Range("L205").Value = 24000
$ws.Range("M205").Value = 24000
$ws.Range("N205").Value = "$/malla 70 unidades"
$ws.Range("O205").Value = "Región de Arica y Parinacota"
$ws.Range("P205").Value = 343
$ws.Range("Q205").Value = 70
$ws.Range("R205").Value = "Hortaliza"
